# Updated cryptos list on Fri Jul  5 09:30:20 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for every coin row, and
# swaps the RenzoRestakedETH / FirstDigitalUSD rows (41/42) back to their
# new ranking order (with their own refreshed price + volume figures).
#
# Note: several "Price" strings look like plain decimal numbers (e.g.
# "465.76", "1.00", "0.0613"). The source workbook stores these as literal
# text, so a leading apostrophe is used for those assignments to force
# Excel to keep them as text (quote-prefixed) instead of silently
# converting them to numeric values and losing the original formatting /
# precision (e.g. "1.00" -> 1, or float rounding noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.242.23"
$ws.Range("E2").Value = "  -5.39%  "
$ws.Range("D3").Value = "2.856.95"
$ws.Range("E3").Value = "  -9.35%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'465.76"
$ws.Range("E5").Value = "  -11.47%  "
$ws.Range("D6").Value = "'124.39"
$ws.Range("E6").Value = "  -7.30%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "2.861.34"
$ws.Range("E8").Value = "  -9.15%  "
$ws.Range("D9").Value = "'0.403"
$ws.Range("E9").Value = "  -10.49%  "
$ws.Range("D10").Value = "'6.57"
$ws.Range("E10").Value = "  -8.97%  "
$ws.Range("D11").Value = "'0.0958"
$ws.Range("E11").Value = "  -13.81%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  -15.40%  "
$ws.Range("E13").Value = "  -4.69%  "
$ws.Range("D14").Value = "3.342.82"
$ws.Range("E14").Value = "  -9.54%  "
$ws.Range("D15").Value = "'23.15"
$ws.Range("E15").Value = "  -10.12%  "
$ws.Range("D16").Value = "54.162.37"
$ws.Range("E16").Value = "  -5.53%  "
$ws.Range("D17").Value = "2.848.78"
$ws.Range("E17").Value = "  -9.72%  "
$ws.Range("E18").Value = "  -13.55%  "
$ws.Range("D19").Value = "'5.33"
$ws.Range("E19").Value = "  -8.38%  "
$ws.Range("D20").Value = "'11.34"
$ws.Range("E20").Value = "  -13.38%  "
$ws.Range("D21").Value = "'7.01"
$ws.Range("E21").Value = "  -12.85%  "
$ws.Range("D22").Value = "'297.08"
$ws.Range("E22").Value = "  -15.04%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'0.435"
$ws.Range("E24").Value = "  -14.88%  "
$ws.Range("D25").Value = "'58.21"
$ws.Range("E25").Value = "  -15.93%  "
$ws.Range("D26").Value = "'0.992"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'0.150"
$ws.Range("E27").Value = "  -9.89%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "0.0₃0788"
$ws.Range("E29").Value = "  -18.53%  "
$ws.Range("D30").Value = "'6.07"
$ws.Range("E30").Value = "  -11.49%  "
$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = "  -11.69%  "
$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = "  -9.06%  "
$ws.Range("D33").Value = "'18.65"
$ws.Range("E33").Value = "  -13.71%  "
$ws.Range("D34").Value = "'1.59"
$ws.Range("E34").Value = "  -15.23%  "
$ws.Range("D35").Value = "'137.69"
$ws.Range("E35").Value = "  -12.90%  "
$ws.Range("D36").Value = "'4.13"
$ws.Range("E36").Value = "  -16.21%  "
$ws.Range("D37").Value = "'5.35"
$ws.Range("E37").Value = "  -14.15%  "
$ws.Range("E38").Value = "  -15.34%  "
$ws.Range("D39").Value = "'22.83"
$ws.Range("E39").Value = "  -10.90%  "
$ws.Range("D40").Value = "'0.0613"
$ws.Range("E40").Value = "  -12.37%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "2.878.92"
$ws.Range("E41").Value = "  -9.38%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'34.93"
$ws.Range("E43").Value = "  -13.06%  "
$ws.Range("E44").Value = "  -13.84%  "
$ws.Range("D45").Value = "'0.927"
$ws.Range("E45").Value = "  -15.18%  "
$ws.Range("D46").Value = "'3.39"
$ws.Range("E46").Value = "  -14.39%  "
$ws.Range("D47").Value = "'1.30"
$ws.Range("E47").Value = "  -11.11%  "
$ws.Range("D48").Value = "2.027.22"
$ws.Range("E48").Value = "  -10.27%  "
$ws.Range("D49").Value = "'5.32"
$ws.Range("E49").Value = "  -14.31%  "
$ws.Range("D50").Value = "'0.0213"
$ws.Range("E50").Value = "  -9.62%  "
$ws.Range("D51").Value = "'17.52"
$ws.Range("E51").Value = "  -14.91%  "
